$d = $word.ActiveDocument

# --- 1) Update the two MERGEFIELD instructions to use the
#        "_in_human_format" fields -------------------------------------
foreach ($f in $d.Fields) {
    $code = $f.Code.Text

    if ($code -like "*MERGEFIELD =occupation_standard.work_processes_hours \*MERGEFORMAT*" -or
        $code -like "*MERGEFIELD =occupation_standard.work_processes_hours \* MERGEFORMAT*") {
        $f.Data = " MERGEFIELD =occupation_standard.work_processes_hours_in_human_format \* MERGEFORMAT "
    }
    elseif ($code -like "*MERGEFIELD =occupation_standard.related_instructions_hours \*MERGEFORMAT*" -or
            $code -like "*MERGEFIELD =occupation_standard.related_instructions_hours \* MERGEFORMAT*") {
        $f.Data = " MERGEFIELD =occupation_standard.related_instructions_hours_in_human_format \* MERGEFORMAT "
    }
}

# --- 2) Drop the stray two-space run that trails the (now renamed)
#        related_instructions_hours field result ------------------------
foreach ($f in $d.Fields) {
    $code = $f.Code.Text
    if ($code -like "*related_instructions_hours_in_human_format*") {
        $tail = $d.Range($f.Result.End + 1, $f.Result.End + 3)
        if ($tail.Text -eq "  ") {
            $tail.Delete()
        }
    }
}
